$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.929
$ws.Range("A8").Value = -21.753
$ws.Range("A10").Value = -21.808
$ws.Range("A12").Value = -21.303
$ws.Range("B12").Value = 6.304
$ws.Range("C12").Value = -11.214
$ws.Range("C13").Value = -12.813
$ws.Range("B15").Value = 5.282999999999999
$ws.Range("B17").Value = 4.836999999999999
$ws.Range("A18").Value = -21.858
$ws.Range("C21").Value = -12.907
$ws.Range("C25").Value = -12.305
$ws.Range("B26").Value = 6.022
$ws.Range("B27").Value = 5.529
$ws.Range("B28").Value = 5.355
$ws.Range("C32").Value = -12.394
$ws.Range("C36").Value = -12.776
$ws.Range("A37").Value = -21.204
$ws.Range("B37").Value = 6.893000000000001
$ws.Range("C38").Value = -12.356
$ws.Range("C41").Value = -12.554
$ws.Range("B47").Value = 5.551
$ws.Range("C52").Value = -11.944
$ws.Range("A55").Value = -22.109
$ws.Range("C59").Value = -12.18
$ws.Range("B65").Value = 5.8
$ws.Range("C67").Value = -11.065
$ws.Range("A68").Value = -21.567
$ws.Range("B73").Value = 6.792
$ws.Range("A77").Value = -21.032
$ws.Range("A78").Value = -20.752
$ws.Range("A81").Value = -21.747
$ws.Range("A82").Value = -21.822
$ws.Range("B84").Value = 5.271000000000001
$ws.Range("C84").Value = -12.517
$ws.Range("B85").Value = 5.324000000000001
$ws.Range("C88").Value = -13.351
$ws.Range("C89").Value = -13.791
$ws.Range("B93").Value = 5.587000000000001
$ws.Range("B95").Value = 6.43
$ws.Range("C95").Value = -11.626
$ws.Range("B98").Value = 6.866
$ws.Range("B99").Value = 5.447
$ws.Range("B101").Value = 6.043000000000001
$ws.Range("C105").Value = -12.753
